# Insert a new data row at row 242 (shifting existing rows 242-323 down to
# 243-324) and populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("242:242").Insert()

$ws.Range("A242").Value = 9
$ws.Range("B242").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C242").Value = 'Metropolitana'
$ws.Range("D242").Value = 44809
$ws.Range("E242").Value = 13
$ws.Range("F242").Value = 100112043
$ws.Range("G242").Value = 'Pepino ensalada'
$ws.Range("H242").Value = 'Sin especificar'
$ws.Range("I242").Value = 'Primera'
$ws.Range("J242").Value = 231
$ws.Range("K242").Value = 24000
$ws.Range("L242").Value = 26000
$ws.Range("M242").Value = 25238
$ws.Range("N242").Value = '$/caja 60 unidades'
$ws.Range("O242").Value = 'Región de Arica y Parinacota'
$ws.Range("P242").Value = 421
$ws.Range("Q242").Value = 60
$ws.Range("R242").Value = 'Hortaliza'
